$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking stats stay stored as text, matching the
# workbook's original convention (values stored as text strings).
$ws.Range("C2:F12").NumberFormat = "@"

$ws.Range("C2").Value = "46"
$ws.Range("D2").Value = "41"
$ws.Range("F2").Value = "2"

$ws.Range("C3").Value = "77"
$ws.Range("D3").Value = "51"
$ws.Range("E3").Value = "7"
$ws.Range("F3").Value = "3"

$ws.Range("C4").Value = "29"
$ws.Range("D4").Value = "27"
$ws.Range("E4").Value = "3"
$ws.Range("F4").Value = "1"

$ws.Range("C5").Value = "27"
$ws.Range("D5").Value = "27"
$ws.Range("E5").Value = "2"
$ws.Range("F5").Value = "1"

$ws.Range("C6").Value = "28"
$ws.Range("D6").Value = "25"
$ws.Range("E6").Value = "4"
$ws.Range("F6").Value = "0"

$ws.Range("C7").Value = "132"
$ws.Range("D7").Value = "69"
$ws.Range("E7").Value = "14"
$ws.Range("F7").Value = "7"

$ws.Range("C8").Value = "15"
$ws.Range("D8").Value = "11"
$ws.Range("E8").Value = "1"
$ws.Range("F8").Value = "1"

$ws.Range("C9").Value = "61"
$ws.Range("D9").Value = "49"
$ws.Range("E9").Value = "1"
$ws.Range("F9").Value = "5"

$ws.Range("C11").Value = "17"
$ws.Range("D11").Value = "19"
$ws.Range("F11").Value = "0"

$ws.Range("C12").Value = "21"
$ws.Range("D12").Value = "19"
